$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3190.6365
$ws.Range("I62").Value = 2671.2856
$ws.Range("K62").Value = 2671.2856
$ws.Range("M62").Value = -2047.2856
$ws.Range("H64").Value = 14096544
$ws.Range("I64").Value = 4080102.5
$ws.Range("K64").Value = 4080102.5
$ws.Range("M64").Value = -4079854.5
$ws.Range("H65").Value = 3190.6365
$ws.Range("I65").Value = 2671.2856
$ws.Range("K65").Value = 13356.428
$ws.Range("M65").Value = -10236.428
$ws.Range("H67").Value = 14096544
$ws.Range("I67").Value = 4080102.5
$ws.Range("K67").Value = 4080102.5
$ws.Range("M67").Value = -4079244.5
$ws.Range("H80").Value = 796.3077
$ws.Range("I80").Value = 585.5
$ws.Range("K80").Value = 1756.5
$ws.Range("M80").Value = -758.5
$ws.Range("H83").Value = 796.3077
$ws.Range("I83").Value = 585.5
$ws.Range("K83").Value = 5269.5
$ws.Range("M83").Value = -277.5
$ws.Range("H132").Value = 8940.712
$ws.Range("I132").Value = 2206.1292
$ws.Range("K132").Value = 6618.3876
$ws.Range("M132").Value = -4088.3876
$ws.Range("H137").Value = 8133712
$ws.Range("J137").Value = 17550352
$ws.Range("L137").Value = 52651056
$ws.Range("N137").Value = -52656156
$ws.Range("H138").Value = 5639.685
$ws.Range("I138").Value = 854.15
$ws.Range("J138").Value = 8454.706
$ws.Range("K138").Value = 2562.45
$ws.Range("L138").Value = 25364.118
$ws.Range("M138").Value = 2577.55
$ws.Range("N138").Value = -35644.118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5983.9585
$ws.Range("I32").Value = 3435.8235
$ws.Range("J32").Value = 12172.286
$ws.Range("K32").Value = 3435.8235
$ws.Range("L32").Value = 12172.286
$ws.Range("M32").Value = -3148.8235
$ws.Range("N32").Value = -12746.286
$ws.Range("H61").Value = 10757.389
$ws.Range("I61").Value = 14070.083
$ws.Range("K61").Value = 14070.083
$ws.Range("M61").Value = -13858.083
$ws.Range("H88").Value = 56188.555
$ws.Range("I88").Value = 739.4
$ws.Range("J88").Value = 125500
$ws.Range("K88").Value = 739.4
$ws.Range("L88").Value = 125500
$ws.Range("M88").Value = -333.4
$ws.Range("N88").Value = -126312
$ws.Range("H91").Value = 56188.555
$ws.Range("I91").Value = 739.4
$ws.Range("J91").Value = 125500
$ws.Range("K91").Value = 739.4
$ws.Range("L91").Value = 125500
$ws.Range("M91").Value = 664.6
$ws.Range("N91").Value = -128308
$ws.Range("H102").Value = 312839.12
$ws.Range("I102").Value = 548910.9
$ws.Range("J102").Value = 2218.3684
$ws.Range("K102").Value = 548910.9
$ws.Range("L102").Value = 2218.3684
$ws.Range("M102").Value = -547288.9
$ws.Range("N102").Value = -5462.368399999999
$ws.Range("H110").Value = 758661.4399999999
$ws.Range("I110").Value = 1201925.9
$ws.Range("J110").Value = 5111.9
$ws.Range("K110").Value = 1201925.9
$ws.Range("L110").Value = 5111.9
$ws.Range("M110").Value = -1199880.9
$ws.Range("N110").Value = -9201.9
$ws.Range("H136").Value = 10757.389
$ws.Range("I136").Value = 14070.083
$ws.Range("K136").Value = 42210.249
$ws.Range("M136").Value = -39660.249
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 979500.9
$ws.Range("I94").Value = 1054693.2
$ws.Range("K94").Value = 1054693.2
$ws.Range("M94").Value = -1054242.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 45929.43
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 63341.2
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 63341.2
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -64589.2
$ws.Range("H65").Value = 45929.43
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 63341.2
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 316706
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -322946
$ws.Range("H99").Value = 7784.0713
$ws.Range("I99").Value = 3795.6
$ws.Range("K99").Value = 3795.6
$ws.Range("M99").Value = -2297.6
$ws.Range("H126").Value = 7784.0713
$ws.Range("I126").Value = 3795.6
$ws.Range("K126").Value = 11386.8
$ws.Range("M126").Value = -8916.799999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4959
$ws.Range("I122").Value = 1749.6666
$ws.Range("J122").Value = 7709.857
$ws.Range("K122").Value = 5248.9998
$ws.Range("L122").Value = 23129.571
$ws.Range("M122").Value = -2798.9998
$ws.Range("N122").Value = -28029.571
$ws.Range("H132").Value = 4486.5264
$ws.Range("I132").Value = 4439.1924
$ws.Range("J132").Value = 4589.0835
$ws.Range("K132").Value = 13317.5772
$ws.Range("L132").Value = 13767.2505
$ws.Range("M132").Value = -10787.5772
$ws.Range("N132").Value = -18827.2505
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7206.6665
$ws.Range("I7").Value = 5682.6665
$ws.Range("J7").Value = 8222.666999999999
$ws.Range("K7").Value = 5682.6665
$ws.Range("L7").Value = 8222.666999999999
$ws.Range("M7").Value = -5570.6665
$ws.Range("N7").Value = -8446.666999999999
$ws.Range("H16").Value = 6251606
$ws.Range("I16").Value = 7408800.5
$ws.Range("J16").Value = 2755.8
$ws.Range("K16").Value = 7408800.5
$ws.Range("L16").Value = 2755.8
$ws.Range("M16").Value = -7408630.5
$ws.Range("N16").Value = -3095.8
$ws.Range("H40").Value = 41669590
$ws.Range("I40").Value = 3099.5
$ws.Range("K40").Value = 3099.5
$ws.Range("M40").Value = -2963.5
$ws.Range("H126").Value = 7206.6665
$ws.Range("I126").Value = 5682.6665
$ws.Range("J126").Value = 8222.666999999999
$ws.Range("K126").Value = 17047.9995
$ws.Range("L126").Value = 24668.001
$ws.Range("M126").Value = -14577.9995
$ws.Range("N126").Value = -29608.001
$ws.Range("H133").Value = 92314
$ws.Range("J133").Value = 92314
$ws.Range("L133").Value = 92314
$ws.Range("N133").Value = -97374
$ws.Range("H136").Value = 2733.875
$ws.Range("I136").Value = 2124.4285
$ws.Range("K136").Value = 6373.2855
$ws.Range("M136").Value = -3823.2855
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 772.4545000000001
$ws.Range("I113").Value = 787.3333
$ws.Range("J113").Value = 705.5
$ws.Range("K113").Value = 2361.9999
$ws.Range("L113").Value = 2116.5
$ws.Range("M113").Value = -191.9998999999998
$ws.Range("N113").Value = -6456.5
$ws.Range("H122").Value = 4205.641
$ws.Range("I122").Value = 3515.8484
$ws.Range("K122").Value = 10547.5452
$ws.Range("M122").Value = -8097.5452
$ws.Range("H126").Value = 1487.4546
$ws.Range("I126").Value = 1246
$ws.Range("J126").Value = 2574
$ws.Range("K126").Value = 3738
$ws.Range("L126").Value = 7722
$ws.Range("M126").Value = -1268
$ws.Range("N126").Value = -12662
$ws.Range("I132").Value = 4275083
$ws.Range("J132").Value = 100003920
$ws.Range("K132").Value = 12825249
$ws.Range("L132").Value = 300011760
$ws.Range("M132").Value = -12822719
$ws.Range("N132").Value = -300016820
$ws.Range("H133").Value = 58777
$ws.Range("J133").Value = 58777
$ws.Range("L133").Value = 58777
$ws.Range("N133").Value = -68897
